$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "Report (CSV) : Removed Special characters"
# The header cell C7 contained a leading-space special character
# (" SETTLEMENT DATE"). Replace it with a clean label "ENTRY DATE",
# keeping the cell's existing formatting (border/number-format/font)
# untouched.
$ws.Range("C7").Value = "ENTRY DATE"
